$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 307. This pushes the existing rows
# 307..312 down to 308..313 (dimension grows from A1:T312 to A1:T313).
$ws.Rows.Item(307).Insert()

# Populate the newly inserted row 307 with the new record.
$ws.Range("A307").Value = 5
$ws.Range("B307").Value = "Macroferia Regional de Talca"
$ws.Range("C307").Value = "Maule"
$ws.Range("D307").Value = 44890
$ws.Range("E307").Value = 7
$ws.Range("F307").Value = "Fruta"
$ws.Range("G307").Value = 100108
$ws.Range("H307").Value = "Tropicales y subtropicales"
$ws.Range("I307").Value = 100108005
$ws.Range("J307").Value = "Piña"
$ws.Range("K307").Value = "Caramelo"
$ws.Range("L307").Value = "Especial"
$ws.Range("M307").Value = 150
$ws.Range("N307").Value = 25000
$ws.Range("O307").Value = 25000
$ws.Range("P307").Value = 25000
$ws.Range("Q307").Value = "$/caja 14 unidades"
$ws.Range("R307").Value = "Ecuador"
$ws.Range("S307").Value = 1786
$ws.Range("T307").Value = 14
